$wb = $excel.ActiveWorkbook

# The localization status changed from "Ready for handoff" to "In Translation"
# on every sheet that reports it (Overview summary columns + the per-locale
# Status column on each language sheet).
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# With the shorter status text, the report generator narrows the status
# columns to fit (AutoFit-style re-generation) on every sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
